$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = 988
$ws.Range("D21").Value = 5608644
$ws.Range("E21").Value = 899.830579175357
$ws.Range("G21").Value = 4.219409282700415
$ws.Range("H21").Value = 28.02282233257824
